# SMP_PP_e&p.xlsx — update the constant "AO" driver value (used by the
# equilibrium-equation check across rows 3:18) from 71011.093371362236 to
# 71636.837037504025. All dependent formulas in columns AR/AS/AT/AU
# recalculate automatically from this single input change.
#
# Commit: "Check validity of equilibrium equations using PE, 0.5PP, PP and
# 2PP cases"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("AO3:AO18")
$rng.Value = 71636.837037504025

# Mirror the author's final selection (AO3 active, AO3:AO18 selected) as
# seen in the saved sheetView of the edited workbook.
$rng.Select()
